$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_12_0_0"
$ws.Range("B2").Value = -0.1358038198910239
$ws.Range("C2").Value = -0.737099494975237
$ws.Range("D2").Value = -0.7385932846899437
$ws.Range("E2").Value = -0.7129251277333457
$ws.Range("F2").Value = 1.257000207901001
$ws.Range("G2").Value = 2.146814584732056
$ws.Range("H2").Value = 1.000275254249573
$ws.Range("I2").Value = 1.607266783714294

$ws.Range("A3").Value = "model_12_0_23"
$ws.Range("B3").Value = 0.03258831475910329
$ws.Range("C3").Value = -1.442879260315176
$ws.Range("D3").Value = -1.425437634459964
$ws.Range("E3").Value = -1.403244618808951
$ws.Range("F3").Value = 1.070639729499817
$ws.Range("G3").Value = 3.019060850143433
$ws.Range("H3").Value = 1.39544141292572
$ws.Range("I3").Value = 2.2550048828125

$ws.Range("A4").Value = "model_12_0_22"
$ws.Range("B4").Value = 0.03360754403340394
$ws.Range("C4").Value = -1.438820858020111
$ws.Range("D4").Value = -1.418631626254077
$ws.Range("E4").Value = -1.39845102364478
$ws.Range("F4").Value = 1.069511651992798
$ws.Range("G4").Value = 3.014045238494873
$ws.Range("H4").Value = 1.391525745391846
$ws.Range("I4").Value = 2.250507354736328

$ws.Range("A5").Value = "model_12_0_21"
$ws.Range("B5").Value = 0.03465496075421703
$ws.Range("C5").Value = -1.434619867369574
$ws.Range("D5").Value = -1.411759726397404
$ws.Range("E5").Value = -1.393538598849687
$ws.Range("F5").Value = 1.068352580070496
$ws.Range("G5").Value = 3.008853197097778
$ws.Range("H5").Value = 1.387572050094604
$ws.Range("I5").Value = 2.245897769927979

$ws.Range("A6").Value = "model_12_0_20"
$ws.Range("B6").Value = 0.03572753993838573
$ws.Range("C6").Value = -1.430288633034875
$ws.Range("D6").Value = -1.404873947111576
$ws.Range("E6").Value = -1.388531748994716
$ws.Range("F6").Value = 1.067165374755859
$ws.Range("G6").Value = 3.003500699996948
$ws.Range("H6").Value = 1.383610367774963
$ws.Range("I6").Value = 2.241199970245361

$ws.Range("A7").Value = "model_12_0_19"
$ws.Range("B7").Value = 0.03682737207193176
$ws.Range("C7").Value = -1.425764242675171
$ws.Range("D7").Value = -1.398480670193544
$ws.Range("E7").Value = -1.383532160546695
$ws.Range("F7").Value = 1.065948247909546
$ws.Range("G7").Value = 2.9979088306427
$ws.Range("H7").Value = 1.379932165145874
$ws.Range("I7").Value = 2.236508846282959

$ws.Range("A8").Value = "model_12_0_18"
$ws.Range("B8").Value = 0.03798159776581833
$ws.Range("C8").Value = -1.420888401619913
$ws.Range("D8").Value = -1.392584865900811
$ws.Range("E8").Value = -1.378431119113659
$ws.Range("F8").Value = 1.06467080116272
$ws.Range("G8").Value = 2.991882801055908
$ws.Range("H8").Value = 1.376539945602417
$ws.Range("I8").Value = 2.231722354888916

$ws.Range("A9").Value = "model_12_0_17"
$ws.Range("B9").Value = 0.03924164540156794
$ws.Range("C9").Value = -1.415584857124467
$ws.Range("D9").Value = -1.385940825140553
$ws.Range("E9").Value = -1.372815970115369
$ws.Range("F9").Value = 1.063276290893555
$ws.Range("G9").Value = 2.985328674316406
$ws.Range("H9").Value = 1.372717499732971
$ws.Range("I9").Value = 2.226453542709351

$ws.Range("A10").Value = "model_12_0_16"
$ws.Range("B10").Value = 0.04063134755311693
$ws.Range("C10").Value = -1.409976356913135
$ws.Range("D10").Value = -1.377671817131423
$ws.Range("E10").Value = -1.36651928160785
$ws.Range("F10").Value = 1.06173837184906
$ws.Range("G10").Value = 2.978397130966187
$ws.Range("H10").Value = 1.367959976196289
$ws.Range("I10").Value = 2.220545053482056

$ws.Range("A11").Value = "model_12_0_15"
$ws.Range("B11").Value = 0.04366143788740329
$ws.Range("C11").Value = -1.398612912476632
$ws.Range("D11").Value = -1.354717508261724
$ws.Range("E11").Value = -1.351972396634954
$ws.Range("F11").Value = 1.058384895324707
$ws.Range("G11").Value = 2.964353561401367
$ws.Range("H11").Value = 1.354753613471985
$ws.Range("I11").Value = 2.206895589828491

$ws.Range("A12").Value = "model_12_0_14"
$ws.Range("B12").Value = 0.04639370565280976
$ws.Range("C12").Value = -1.387513316717293
$ws.Range("D12").Value = -1.337676781964104
$ws.Range("E12").Value = -1.339315413458546
$ws.Range("F12").Value = 1.055361151695251
$ws.Range("G12").Value = 2.950636148452759
$ws.Range("H12").Value = 1.34494948387146
$ws.Range("I12").Value = 2.195019245147705

$ws.Range("A13").Value = "model_12_0_13"
$ws.Range("B13").Value = 0.05080809512079898
$ws.Range("C13").Value = -1.369091374908458
$ws.Range("D13").Value = -1.311577013759429
$ws.Range("E13").Value = -1.318938841084485
$ws.Range("F13").Value = 1.050475835800171
$ws.Range("G13").Value = 2.927869081497192
$ws.Range("H13").Value = 1.329933285713196
$ws.Range("I13").Value = 2.175899505615234

$ws.Range("A14").Value = "model_12_0_12"
$ws.Range("B14").Value = 0.05765007860785487
$ws.Range("C14").Value = -1.339127345614866
$ws.Range("D14").Value = -1.274356659165706
$ws.Range("E14").Value = -1.287305590334351
$ws.Range("F14").Value = 1.042903661727905
$ws.Range("G14").Value = 2.890837669372559
$ws.Range("H14").Value = 1.308519124984741
$ws.Range("I14").Value = 2.146217584609985

$ws.Range("A15").Value = "model_12_0_11"
$ws.Range("B15").Value = 0.06455591439513142
$ws.Range("C15").Value = -1.309961250715151
$ws.Range("D15").Value = -1.231743120010019
$ws.Range("E15").Value = -1.254672488085743
$ws.Range("F15").Value = 1.035260915756226
$ws.Range("G15").Value = 2.854792594909668
$ws.Range("H15").Value = 1.284002065658569
$ws.Range("I15").Value = 2.115597248077393

$ws.Range("A16").Value = "model_12_0_10"
$ws.Range("B16").Value = 0.07332995286323918
$ws.Range("C16").Value = -1.270552062811575
$ws.Range("D16").Value = -1.18629818216301
$ws.Range("E16").Value = -1.214079919300428
$ws.Range("F16").Value = 1.025550603866577
$ws.Range("G16").Value = 2.806087970733643
$ws.Range("H16").Value = 1.257855892181396
$ws.Range("I16").Value = 2.077508926391602

$ws.Range("A17").Value = "model_12_0_9"
$ws.Range("B17").Value = 0.07671670308132283
$ws.Range("C17").Value = -1.253584371378263
$ws.Range("D17").Value = -1.153693389460128
$ws.Range("E17").Value = -1.192840391670041
$ws.Range("F17").Value = 1.021802544593811
$ws.Range("G17").Value = 2.785118579864502
$ws.Range("H17").Value = 1.239097118377686
$ws.Range("I17").Value = 2.057579278945923

$ws.Range("A18").Value = "model_12_0_2"
$ws.Range("B18").Value = 0.07988313780032585
$ws.Range("C18").Value = -1.061156548177444
$ws.Range("D18").Value = -0.9765742014728127
$ws.Range("E18").Value = -1.007555474918414
$ws.Range("F18").Value = 1.018298268318176
$ws.Range("G18").Value = 2.547304391860962
$ws.Range("H18").Value = 1.137194156646729
$ws.Range("I18").Value = 1.883723258972168

$ws.Range("A19").Value = "model_12_0_8"
$ws.Range("B19").Value = 0.08052912847130733
$ws.Range("C19").Value = -1.235445075250986
$ws.Range("D19").Value = -1.138884225921873
$ws.Range("E19").Value = -1.175919031062208
$ws.Range("F19").Value = 1.017583250999451
$ws.Range("G19").Value = 2.762700796127319
$ws.Range("H19").Value = 1.230576992034912
$ws.Range("I19").Value = 2.041701793670654

$ws.Range("A20").Value = "model_12_0_7"
$ws.Range("B20").Value = 0.08119871015771307
$ws.Range("C20").Value = -1.228424142079101
$ws.Range("D20").Value = -1.13877244505624
$ws.Range("E20").Value = -1.170991244401677
$ws.Range("F20").Value = 1.016842246055603
$ws.Range("G20").Value = 2.754023790359497
$ws.Range("H20").Value = 1.230512619018555
$ws.Range("I20").Value = 2.037077903747559

$ws.Range("A21").Value = "model_12_0_6"
$ws.Range("B21").Value = 0.08683375296157547
$ws.Range("C21").Value = -1.194309627843353
$ws.Range("D21").Value = -1.111964542785348
$ws.Range("E21").Value = -1.139468299152865
$ws.Range("F21").Value = 1.010605931282043
$ws.Range("G21").Value = 2.711863279342651
$ws.Range("H21").Value = 1.215089082717896
$ws.Range("I21").Value = 2.007499217987061

$ws.Range("A22").Value = "model_12_0_5"
$ws.Range("B22").Value = 0.08692879707774848
$ws.Range("C22").Value = -1.142556205037639
$ws.Range("D22").Value = -1.153878788249778
$ws.Range("E22").Value = -1.115474937483443
$ws.Range("F22").Value = 1.01050078868866
$ws.Range("G22").Value = 2.647902965545654
$ws.Range("H22").Value = 1.239203929901123
$ws.Range("I22").Value = 1.984986066818237

$ws.Range("A23").Value = "model_12_0_4"
$ws.Range("B23").Value = 0.0884316081703963
$ws.Range("C23").Value = -1.127992885851564
$ws.Range("D23").Value = -1.13477582919973
$ws.Range("E23").Value = -1.099808123254763
$ws.Range("F23").Value = 1.008837699890137
$ws.Range("G23").Value = 2.629904747009277
$ws.Range("H23").Value = 1.228213310241699
$ws.Range("I23").Value = 1.970285534858704

$ws.Range("A24").Value = "model_12_0_3"
$ws.Range("B24").Value = 0.1265698161120716
$ws.Range("C24").Value = -0.9919137749496807
$ws.Range("D24").Value = -0.8198425067634416
$ws.Range("E24").Value = -0.9140489823196181
$ws.Range("F24").Value = 0.9666298627853394
$ws.Range("G24").Value = 2.461730003356934
$ws.Range("H24").Value = 1.047020792961121
$ws.Range("I24").Value = 1.795984625816345

$ws.Range("A25").Value = "model_12_0_1"
$ws.Range("B25").Value = 0.174629113481016
$ws.Range("C25").Value = -0.6498254584078758
$ws.Range("D25").Value = -0.4100896396978513
$ws.Range("E25").Value = -0.5572818313510248
$ws.Range("F25").Value = 0.9134424328804016
$ws.Range("G25").Value = 2.038956165313721
$ws.Range("H25").Value = 0.8112752437591553
$ws.Range("I25").Value = 1.461223959922791

$ws.Rows.Item(26).Delete()

Write-Output "done"